$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 98
$ws.Range("H2").Value = 30
$ws.Range("I2").Value = 38
$ws.Range("J2").Value = 166

# Row 3 updates
$ws.Range("H3").Value = 25
$ws.Range("J3").Value = 155
